# ---------------------------------------------------------------------
# Acta-reunión-Seguimiento2.docx - "cambio"
#
# 1) Drop the _GoBack bookmark that currently sits right after
#    "...desviación de avance de " (it is being relocated elsewhere
#    in the document).
# 2) Change "Próxima Reunión: 02 de octubre de 2015" to
#    "... : 03 ..." and re-insert the _GoBack bookmark right after
#    the new "03", splitting the old ": 02 " run into ": 03" plus a
#    standalone run that holds the trailing space.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: remove the old _GoBack bookmark -------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Change 2: "02" -> "03" plus a relocated _GoBack bookmark ----------

# Find the paragraph that holds "Próxima Reunión: 02 ..."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains(": 02 ")) {
        $target = $cand
        break
    }
}

$pBase = $target.Range.Start
$pText = $target.Range.Text
$marker = ": 02 "
$markerPos = $pText.IndexOf($marker)

# Position right after "02" (before the trailing space) - this is
# where the _GoBack bookmark needs to end up.
$afterDigits = $pBase + $markerPos + 4
# Position right after the trailing space (before the following "de ")
# - used only temporarily, so the text edit below cannot bleed past
# the space and swallow the untouched runs that follow it.
$afterSpace = $pBase + $markerPos + 5

$d.Bookmarks.Add("_GoBack", $d.Range($afterDigits, $afterDigits))
$d.Bookmarks.Add("zzTempSplit", $d.Range($afterSpace, $afterSpace))

# Replace the digits "02" with "03".
$pText2 = $target.Range.Text
$digitsStart = $pBase + $pText2.IndexOf("02")
$digitsEnd = $digitsStart + 2
$d.Range($digitsStart, $digitsEnd).Text = "03"

# Re-touch the lone space run sitting between the two bookmarks, so it
# ends up as its own clean run (mirroring what happened on the "02"
# run above) instead of keeping stale formatting-run metadata.
$pText3 = $target.Range.Text
$spaceStart = $pBase + ($pText3.IndexOf(": 03") + 4)
$spaceEnd = $spaceStart + 1
$d.Range($spaceStart, $spaceEnd).Text = "#"
$pText4 = $target.Range.Text
$hashStart = $pBase + $pText4.IndexOf("#")
$hashEnd = $hashStart + 1
$d.Range($hashStart, $hashEnd).Text = " "

# Drop the temporary helper bookmark - its only job was to keep the
# edit above from merging into the untouched "de octubre de 2015" runs.
$d.Bookmarks.Item("zzTempSplit").Delete()
